$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Objetivos:" (row 10) B/C content: was wrongly showing the
# "Docentes responsaveis" teacher-name text; correct Portuguese objectives. ---
$ws.Range("B10").Value = "Proporcionar aos discentes os conhecimentos de biologia celular necessários à compreensão das demais disciplinas do curso e a formação do Engenheiro Ambiental."
$ws.Range("C10").Value = "Proporcionar aos discentes os conhecimentos de biologia celular necessários à compreensão das demais disciplinas do curso e a formação do Engenheiro Ambiental."

# --- Insert a new row 13 ("Docentes responsaveis:" content row) so that every
# label in column A lines up with its correct B/C content on the same row
# (everything below was previously shifted up by one row). ---
$ws.Rows("13:13").Insert()
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "1304060 - Maria das Graças de Almeida Felipe"
$ws.Range("C13").Value = "1304060 - Maria das Graças de Almeida Felipe"

# --- Row 14: "Programa resumido:" now gets the correct Portuguese short
# syllabus text (previously held "Semestral"). ---
$ws.Range("B14").Value = "Análise estrutural das células ao microscópio; moléculas orgânicas; organização interna da célula; organelas celulares transdutoras de energia; material genético e mecanismo de divisão celular."
$ws.Range("C14").Value = "Análise estrutural das células ao microscópio; moléculas orgânicas; organização interna da célula; organelas celulares transdutoras de energia; material genético e mecanismo de divisão celular."

# --- Row 16: "Programa:" now gets the correct Portuguese full syllabus text
# (previously held the activation date "01/01/2022"). ---
$ws.Range("B16").Value = "- Estrutura celular e história evolutiva: microrganismos procarióticos eeucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea eEukarya.- Análise estrutural das células ao microscópio: microscopia ótica e eletrônica.- Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos, ácidos nucleicos, aminoácidos. - Organização interna da célula: estrutura e função da membrana plasmática; compartimentos intracelulares e seleção de proteínas; tráfico de vesículas (via de exocitose e endocitose).- Núcleo e organização do material genético: estrutura e função- Ciclo celular e divisão celular: mitose e meiose.- Organelas celulares transdutoras de energia: mitocôndria e cloroplasto."
$ws.Range("C16").Value = "- Estrutura celular e história evolutiva: microrganismos procarióticos eeucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea eEukarya.- Análise estrutural das células ao microscópio: microscopia ótica e eletrônica.- Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos, ácidos nucleicos, aminoácidos. - Organização interna da célula: estrutura e função da membrana plasmática; compartimentos intracelulares e seleção de proteínas; tráfico de vesículas (via de exocitose e endocitose).- Núcleo e organização do material genético: estrutura e função- Ciclo celular e divisão celular: mitose e meiose.- Organelas celulares transdutoras de energia: mitocôndria e cloroplasto."

# --- Row 19: "Método:" now gets the correct evaluation-method text
# (previously held the "1304060 - Maria..." teacher-name text). ---
$ws.Range("B19").Value = "Duas provas escritas (P1 e P2) distribuídas no semestre.Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."
$ws.Range("C19").Value = "Duas provas escritas (P1 e P2) distribuídas no semestre.Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."

# --- Row 20: "Critério:" now gets the correct MF-formula text (previously
# held the "Duas provas escritas..." method text). ---
$ws.Range("B20").Value = "MF = média finalMF = (P1 + P2)/2"
$ws.Range("C20").Value = "MF = média finalMF = (P1 + P2)/2"

# --- Row 21: "Norma de recuperação:" now gets the correct recovery-rule text
# (previously held the "MF = média final..." formula text). ---
$ws.Range("B21").Value = "Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."
$ws.Range("C21").Value = "Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0."

# --- Row 22 (new, previously didn't exist): "Bibliografia:" now gets its own
# row with label + B/C content (previously the "Nota final..." text sat one
# row up, directly under "Bibliografia:" with no row of its own). ---
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "-Alberts, B. et al. Biologia Molecular da Célula, 5ed. Artmed Editora Ltda, 2010.-Cooper, G.M.; Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3ª Edição, 2007.-Wasserman, S.A.; Monorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora, 8ª Edição, 2010.-Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.-Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14 Edição, 2016. -De Roberts, E.M.F.; Hibs, J. Bases da biologia celular e molecular. Editora Guanabara Koogan, 2006.-Taiz, L.; Zeiger, E. Plant Physiology. Mass. Sinauer Associates, 2006."
$ws.Range("C22").Value = "-Alberts, B. et al. Biologia Molecular da Célula, 5ed. Artmed Editora Ltda, 2010.-Cooper, G.M.; Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3ª Edição, 2007.-Wasserman, S.A.; Monorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora, 8ª Edição, 2010.-Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.-Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14 Edição, 2016. -De Roberts, E.M.F.; Hibs, J. Bases da biologia celular e molecular. Editora Guanabara Koogan, 2006.-Taiz, L.; Zeiger, E. Plant Physiology. Mass. Sinauer Associates, 2006."

Write-Host "Edit complete"
